$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("G2").Value = 3.9
$ws.Range("H2").Value = 3.9
$ws.Range("I2").Value = 1.8
$ws.Range("J2").Value = 4
$ws.Range("K2").Value = 2.4
$ws.Range("L2").Value = 2.38
$ws.Range("O2").Value = 1.17
$ws.Range("P2").Value = 5
$ws.Range("Q2").Value = 1.53
$ws.Range("R2").Value = 2.4
$ws.Range("S2").Value = 1.29
$ws.Range("T2").Value = 3.5
$ws.Range("W2").Value = 15
$ws.Range("X2").Value = 23
$ws.Range("Y2").Value = 13
$ws.Range("Z2").Value = 41
$ws.Range("AA2").Value = 26
$ws.Range("AC2").Value = 17
$ws.Range("AD2").Value = 8
$ws.Range("AF2").Value = 34
$ws.Range("AI2").Value = 11
$ws.Range("AK2").Value = 17
$ws.Range("AL2").Value = 13
$ws.Range("AN2").Value = 6
$ws.Range("AO2").Value = 19
$ws.Range("AQ2").Value = 51
$ws.Range("AT2").Value = 3.5
$ws.Range("AU2").Value = 7
$ws.Range("AX2").Value = 4.33
$ws.Range("AY2").Value = 9.5
$ws.Range("AZ2").Value = 17
$ws.Range("BA2").Value = 29

# Row 3
$ws.Range("O3").Value = 1.13
$ws.Range("P3").Value = 6
$ws.Range("Q3").Value = 1.44
$ws.Range("R3").Value = 2.7

# Row 4
$ws.Range("G4").Value = 1.9
$ws.Range("I4").Value = 3.8
$ws.Range("J4").Value = 2.5
$ws.Range("K4").Value = 2.2
$ws.Range("L4").Value = 4.33
$ws.Range("M4").Value = 1.05
$ws.Range("N4").Value = 11
$ws.Range("Q4").Value = 1.85
$ws.Range("R4").Value = 1.95
$ws.Range("U4").Value = 1.73
$ws.Range("V4").Value = 2
$ws.Range("W4").Value = 7.5
$ws.Range("X4").Value = 9.5
$ws.Range("Y4").Value = 8.5
$ws.Range("AB4").Value = 26
$ws.Range("AC4").Value = 11
$ws.Range("AE4").Value = 15
$ws.Range("AF4").Value = 51
$ws.Range("AG4").Value = 201
$ws.Range("AH4").Value = 11
$ws.Range("AJ4").Value = 13
$ws.Range("AL4").Value = 29
$ws.Range("AO4").Value = 10
$ws.Range("AP4").Value = 21
$ws.Range("AS4").Value = 151
$ws.Range("AU4").Value = 8
$ws.Range("AX4").Value = 6
$ws.Range("AY4").Value = 21
$ws.Range("AZ4").Value = 29
$ws.Range("BA4").Value = 67
$ws.Range("BC4").Value = 201

# Row 5
$ws.Range("G5").Value = 1.7
$ws.Range("H5").Value = 3.4
$ws.Range("I5").Value = 5.25
$ws.Range("J5").Value = 2.5
$ws.Range("L5").Value = 6
$ws.Range("O5").Value = 1.53
$ws.Range("P5").Value = 2.5
$ws.Range("Z5").Value = 13
$ws.Range("AA5").Value = 19
$ws.Range("AC5").Value = 6.5
$ws.Range("AH5").Value = 9.5
$ws.Range("AI5").Value = 23
$ws.Range("AM5").Value = 51

# Row 6
$ws.Range("AW6").Value = 151
$ws.Range("BD6").Value = 151

# Row 7
$ws.Range("O7").Value = 1.25
$ws.Range("P7").Value = 4

# Row 8
$ws.Range("AA8").Value = 26

# Row 9
$ws.Range("Q9").Value = 2.1
$ws.Range("R9").Value = 1.73

# Row 11
$ws.Range("G11").Value = 1.65
$ws.Range("H11").Value = 3.75
$ws.Range("I11").Value = 5
$ws.Range("K11").Value = 2.25
$ws.Range("M11").Value = 1.05
$ws.Range("N11").Value = 11
$ws.Range("Q11").Value = 1.83
$ws.Range("R11").Value = 2.03
$ws.Range("U11").Value = 1.83
$ws.Range("V11").Value = 1.83
$ws.Range("W11").Value = 7
$ws.Range("X11").Value = 8
$ws.Range("AB11").Value = 26
$ws.Range("AC11").Value = 11
$ws.Range("AD11").Value = 7
$ws.Range("AG11").Value = 251
$ws.Range("AJ11").Value = 17
$ws.Range("AN11").Value = 3.6
$ws.Range("AS11").Value = 151
$ws.Range("AU11").Value = 8.5
$ws.Range("AZ11").Value = 34
$ws.Range("BA11").Value = 101
$ws.Range("BC11").Value = 251

# Row 14
$ws.Range("G14").Value = 2.27
$ws.Range("I14").Value = 2.92
$ws.Range("J14").Value = 2.82
$ws.Range("L14").Value = 3.5
$ws.Range("S14").Value = 1.35
$ws.Range("U14").Value = 1.53
$ws.Range("W14").Value = 10
$ws.Range("Y14").Value = 8.75
$ws.Range("Z14").Value = 26
$ws.Range("AA14").Value = 17
$ws.Range("AB14").Value = 21
$ws.Range("AE14").Value = 11
$ws.Range("AI14").Value = 17.5
$ws.Range("AJ14").Value = 10.25
$ws.Range("AK14").Value = 40
$ws.Range("AL14").Value = 24
$ws.Range("AM14").Value = 26
$ws.Range("AN14").Value = 4.45
$ws.Range("AO14").Value = 12
$ws.Range("AP14").Value = 16.5
$ws.Range("AQ14").Value = 45
$ws.Range("AR14").Value = 65
$ws.Range("AS14").Value = 175
$ws.Range("AU14").Value = 6.2
$ws.Range("AX14").Value = 5.2
$ws.Range("AY14").Value = 16.5
$ws.Range("AZ14").Value = 20
$ws.Range("BA14").Value = 75
$ws.Range("BB14").Value = 90
